# Edgar_scaling_mapping.xlsx debugging updates
#
# 1) Add a new scaling-override row on the "year" sheet so that Iran (irn)
#    1A3e emissions are only scaled to Edgar from 1992-2011 (not through
#    2012 like the blanket rule) - the 2012 scaling factor is too large and
#    creates an unrealistic jump in "other transport" emissions.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("year")

$ws.Range("A33").Value = "irn"
$ws.Range("B33").Value = "1A3e"
$ws.Range("C33").Value = "NA"
$ws.Range("D33").Value = "NA"
$ws.Range("E33").Value = "NA"
$ws.Range("F33").Value = 1992
$ws.Range("G33").Value = 2011
$ws.Range("H33").Value = "Avoid jump in other tranpsort emissions in 2012 (scaling factor extended and is too large for years past 2012)"

# Restore the on-screen selections that Excel records in the sheet views:
# the "map" sheet had its selected cell scrolled up to C23 ...
$wsMap = $wb.Worksheets.Item("map")
$wsMap.Activate()
$wsMap.Range("C23").Select()

# ... and the "year" sheet (the active tab) ends up with F34 selected,
# just below the newly added row.
$ws.Activate()
$ws.Range("F34").Select()
